# Apply the NaT_inf_inact / NaT_inf_act additions to iAMC_kinetics.xlsx
#
# Adds two new shared strings ("Vmembrane", "G_Gmax_mean"), a new "raw
# 0.00E+00" number-format style used by the SD row's first data cell on the
# NaT_inf_inact sheet, a bestFit column A, and four new data rows (6-9:
# Vmembrane, G_Gmax_mean, SD, SEM) on both the NaT_inf_inact and
# NaT_inf_act worksheets.

$wb = $excel.ActiveWorkbook

$ws5 = $wb.Worksheets.Item("NaT_inf_inact")
$ws6 = $wb.Worksheets.Item("NaT_inf_act")

# --- Row labels (column A) -------------------------------------------------
$ws5.Cells.Item(6, 1).Value = "Vmembrane"
$ws5.Cells.Item(7, 1).Value = "G_Gmax_mean"
$ws5.Cells.Item(8, 1).Value = "SD"
$ws5.Cells.Item(9, 1).Value = "SEM"

$ws6.Cells.Item(6, 1).Value = "Vmembrane"
$ws6.Cells.Item(7, 1).Value = "G_Gmax_mean"
$ws6.Cells.Item(8, 1).Value = "SD"
$ws6.Cells.Item(9, 1).Value = "SEM"

# --- Data rows (columns B:Q) ------------------------------------------------

# NaT_inf_inact (sheet5) row 6: Vmembrane
$ws5.Cells.Item(6, 2).Value = -100
$ws5.Cells.Item(6, 3).Value = -90
$ws5.Cells.Item(6, 4).Value = -80
$ws5.Cells.Item(6, 5).Value = -70
$ws5.Cells.Item(6, 6).Value = -60
$ws5.Cells.Item(6, 7).Value = -50
$ws5.Cells.Item(6, 8).Value = -40
$ws5.Cells.Item(6, 9).Value = -30
$ws5.Cells.Item(6, 10).Value = -20
$ws5.Cells.Item(6, 11).Value = -10
$ws5.Cells.Item(6, 12).Value = 0
$ws5.Cells.Item(6, 13).Value = 10
$ws5.Cells.Item(6, 14).Value = 20
$ws5.Cells.Item(6, 15).Value = 30
$ws5.Cells.Item(6, 16).Value = 40
$ws5.Cells.Item(6, 17).Value = 50

# NaT_inf_inact (sheet5) row 7: G_Gmax_mean
$ws5.Cells.Item(7, 2).Value = 0.93
$ws5.Cells.Item(7, 3).Value = 0.89
$ws5.Cells.Item(7, 4).Value = 0.98
$ws5.Cells.Item(7, 5).Value = 0.99
$ws5.Cells.Item(7, 6).Value = 1
$ws5.Cells.Item(7, 7).Value = 0.822734583
$ws5.Cells.Item(7, 8).Value = 0.460393948
$ws5.Cells.Item(7, 9).Value = 0.223825974
$ws5.Cells.Item(7, 10).Value = 0.114555137
$ws5.Cells.Item(7, 11).Value = 0.058905299
$ws5.Cells.Item(7, 12).Value = 0
$ws5.Cells.Item(7, 13).Value = 0
$ws5.Cells.Item(7, 14).Value = 0
$ws5.Cells.Item(7, 15).Value = 0
$ws5.Cells.Item(7, 16).Value = 0
$ws5.Cells.Item(7, 17).Value = 0

# NaT_inf_inact (sheet5) row 8: SD
$ws5.Cells.Item(8, 2).Value = 0
$ws5.Cells.Item(8, 3).Value = 0.00281692
$ws5.Cells.Item(8, 4).Value = 0.0151037
$ws5.Cells.Item(8, 5).Value = 0.00908173
$ws5.Cells.Item(8, 6).Value = 0.0293016
$ws5.Cells.Item(8, 7).Value = 0.118613
$ws5.Cells.Item(8, 8).Value = 0.159719
$ws5.Cells.Item(8, 9).Value = 0.127634
$ws5.Cells.Item(8, 10).Value = 0.086147
$ws5.Cells.Item(8, 11).Value = 0.0652507
$ws5.Cells.Item(8, 12).Value = 0
$ws5.Cells.Item(8, 13).Value = 0
$ws5.Cells.Item(8, 14).Value = 0
$ws5.Cells.Item(8, 15).Value = 0
$ws5.Cells.Item(8, 16).Value = 0
$ws5.Cells.Item(8, 17).Value = 0

# NaT_inf_inact (sheet5) row 9: SEM
$ws5.Cells.Item(9, 2).Value = 0
$ws5.Cells.Item(9, 3).Value = 0.00199186
$ws5.Cells.Item(9, 4).Value = 0.01068
$ws5.Cells.Item(9, 5).Value = 0.00642175
$ws5.Cells.Item(9, 6).Value = 0.00976721
$ws5.Cells.Item(9, 7).Value = 0.0395376
$ws5.Cells.Item(9, 8).Value = 0.0532397
$ws5.Cells.Item(9, 9).Value = 0.0425446
$ws5.Cells.Item(9, 10).Value = 0.0287157
$ws5.Cells.Item(9, 11).Value = 0.0217502
$ws5.Cells.Item(9, 12).Value = 0
$ws5.Cells.Item(9, 13).Value = 0
$ws5.Cells.Item(9, 14).Value = 0
$ws5.Cells.Item(9, 15).Value = 0
$ws5.Cells.Item(9, 16).Value = 0
$ws5.Cells.Item(9, 17).Value = 0

# NaT_inf_act (sheet6) row 6: Vmembrane
$ws6.Cells.Item(6, 2).Value = -100
$ws6.Cells.Item(6, 3).Value = -90
$ws6.Cells.Item(6, 4).Value = -80
$ws6.Cells.Item(6, 5).Value = -70
$ws6.Cells.Item(6, 6).Value = -60
$ws6.Cells.Item(6, 7).Value = -50
$ws6.Cells.Item(6, 8).Value = -40
$ws6.Cells.Item(6, 9).Value = -30
$ws6.Cells.Item(6, 10).Value = -20
$ws6.Cells.Item(6, 11).Value = -10
$ws6.Cells.Item(6, 12).Value = 0
$ws6.Cells.Item(6, 13).Value = 10
$ws6.Cells.Item(6, 14).Value = 20
$ws6.Cells.Item(6, 15).Value = 30
$ws6.Cells.Item(6, 16).Value = 40
$ws6.Cells.Item(6, 17).Value = 50

# NaT_inf_act (sheet6) row 7: G_Gmax_mean
$ws6.Cells.Item(7, 2).Value = 0
$ws6.Cells.Item(7, 3).Value = 0
$ws6.Cells.Item(7, 4).Value = 0
$ws6.Cells.Item(7, 5).Value = 0
$ws6.Cells.Item(7, 6).Value = 0
$ws6.Cells.Item(7, 7).Value = 0.151134072
$ws6.Cells.Item(7, 8).Value = 0.652125652
$ws6.Cells.Item(7, 9).Value = 0.710218052
$ws6.Cells.Item(7, 10).Value = 0.729662792
$ws6.Cells.Item(7, 11).Value = 0.754033648
$ws6.Cells.Item(7, 12).Value = 0.785162624
$ws6.Cells.Item(7, 13).Value = 0.808873612
$ws6.Cells.Item(7, 14).Value = 0.842014159
$ws6.Cells.Item(7, 15).Value = 0.887992959
$ws6.Cells.Item(7, 16).Value = 0.940687911
$ws6.Cells.Item(7, 17).Value = 1

# NaT_inf_act (sheet6) row 8: SD
$ws6.Cells.Item(8, 2).Value = 0
$ws6.Cells.Item(8, 3).Value = 0
$ws6.Cells.Item(8, 4).Value = 0
$ws6.Cells.Item(8, 5).Value = 0
$ws6.Cells.Item(8, 6).Value = 0
$ws6.Cells.Item(8, 7).Value = 0.436659
$ws6.Cells.Item(8, 8).Value = 0.0176435
$ws6.Cells.Item(8, 9).Value = 0.0300187
$ws6.Cells.Item(8, 10).Value = 0.0383572
$ws6.Cells.Item(8, 11).Value = 0.0555955
$ws6.Cells.Item(8, 12).Value = 0.0677663
$ws6.Cells.Item(8, 13).Value = 0.0862523
$ws6.Cells.Item(8, 14).Value = 0.107416
$ws6.Cells.Item(8, 15).Value = 0.127455
$ws6.Cells.Item(8, 16).Value = 0.153425
$ws6.Cells.Item(8, 17).Value = 0.175676

# NaT_inf_act (sheet6) row 9: SEM
$ws6.Cells.Item(9, 2).Value = 0
$ws6.Cells.Item(9, 3).Value = 0
$ws6.Cells.Item(9, 4).Value = 0
$ws6.Cells.Item(9, 5).Value = 0
$ws6.Cells.Item(9, 6).Value = 0
$ws6.Cells.Item(9, 7).Value = 0.145553
$ws6.Cells.Item(9, 8).Value = 0.00588117
$ws6.Cells.Item(9, 9).Value = 0.0100062
$ws6.Cells.Item(9, 10).Value = 0.0127857
$ws6.Cells.Item(9, 11).Value = 0.0185318
$ws6.Cells.Item(9, 12).Value = 0.0225888
$ws6.Cells.Item(9, 13).Value = 0.0287508
$ws6.Cells.Item(9, 14).Value = 0.0358053
$ws6.Cells.Item(9, 15).Value = 0.0424849
$ws6.Cells.Item(9, 16).Value = 0.0511415
$ws6.Cells.Item(9, 17).Value = 0.0585587

# --- Number format for the new SD row's leading zero cell ------------------
# (mirrors the pre-existing B3/B4 "center, no format" style but without the
#  center alignment - a plain scientific 0.00E+00 number format)
$ws5.Cells.Item(8, 2).NumberFormat = "0.00E+00"
$ws5.Cells.Item(9, 2).NumberFormat = "0.00E+00"

# --- Column A width (auto-fit to the new longer labels) --------------------
$ws5.Columns.Item(1).ColumnWidth = 14.5
$ws6.Columns.Item(1).ColumnWidth = 14.5

# --- Selection / active cell -----------------------------------------------
# NaT_inf_act is selected first so that NaT_inf_inact ends up as the
# workbook's final active sheet/tab, matching the target state.
$ws6.Range("A7").Select()
$ws5.Range("A7").Select()
